# feat: add 2022-Q3 data
#
# Insert a new "2022-Q3" sheet (with its fund-holding detail rows) right
# after "2022-Q2" tab position (i.e. right after the "总计" sheet), and add
# the corresponding summary row to the "总计" sheet. No other sheet's data
# changes - they simply shift one tab position to the right.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet by duplicating the "2022-Q2"
#    sheet (so header row / column styles / borders are carried over
#    faithfully), positioned immediately before "2022-Q2".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("总计").Next()
$q3.Name = "2022-Q3"

# 2022-Q2 had 3 fund rows (rows 2-4); 2022-Q3 only has 2, so drop row 4.
$q3.Rows.Item(4).Delete()

# Row 2: 164811 / 工银瑞信中证京津冀协同发展主题指数（LOF）A (code/name unchanged)
$q3.Range("D2").Value = "'0.12"
$q3.Range("E2").Value = "'93.09"
$q3.Range("F2").Value = "'3.34"
$q3.Range("G2").Value = "'0.0040"
$q3.Range("H2").Value = 2
$q3.Range("D2:G2").Style = "Normal"

# Row 3: 164825 / 工银瑞信中证京津冀协同发展主题指数（LOF）C
$q3.Range("B3").Value = "'164825"
$q3.Range("C3").Value = "工银瑞信中证京津冀协同发展主题指数（LOF）C"
$q3.Range("D3").Value = "'0.03"
$q3.Range("E3").Value = "'93.09"
$q3.Range("F3").Value = "'3.34"
$q3.Range("G3").Value = "'0.0010"
$q3.Range("H3").Value = 2
$q3.Range("B3").Style = "Normal"
$q3.Range("D3:G3").Style = "Normal"

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row right under the
#    header for 2022-Q3, pushing the existing quarters down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Copy index-column formatting down from the row that got pushed to A3.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0
$total.Range("B2:D2").Style = "Normal"

# Re-number the index column (A) for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
